$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.965.98"
$ws.Range("D3").Value = "3.251.66"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.21%  "
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "3.250.28"
$ws.Range("E9").Value = "  +2.44%  "
$ws.Range("E10").Value = "  +4.34%  "
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.408"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.90%  "
$ws.Range("D13").Value = "3.815.57"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("E14").Value = "  +1.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("D16").Value = "66.945.38"
$ws.Range("E16").Value = "  +4.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000167"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.19%  "
$ws.Range("D18").Value = "3.250.46"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("E19").Value = "  +2.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "369.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.61%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.507"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").Value = "3.385.87"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  +3.37%  "
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  +4.98%  "
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "172.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.56%  "
$ws.Range("E36").Value = "  +2.76%  "
$ws.Range("E38").Value = "  +4.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.854"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.31%  "
$ws.Range("E40").Value = "  +10.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.29%  "
$ws.Range("E42").Value = "  +2.03%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.720.60"
$ws.Range("E43").Value = "  +2.49%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.72%  "
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.20%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0673"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.16%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "336.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.96%  "
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("E51").Value = "  +2.57%  "
